$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000000000001257882686900301956
$ws.Range("C2").Value = 0.000002220651329265522090883030
$ws.Range("D2").Value = 0.152905782018181196635353558122
$ws.Range("E2").Value = 6.481428077270620313754534436157
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 6.634336079941388852887484972598
